# "improve handle of empty multiform file"
#
# Adds a new worksheet "2013" (a copy-in-spirit of the existing empty
# "2012" multiform sheet) at the end of the workbook, makes it the active
# sheet/tab, and leaves the previously-active sheet ("2012") no longer
# marked as the selected tab.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the current last sheet so it lands at
# the end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2013"

# Same single "no instruction" placeholder cell the other empty multiform
# sheet ("2012") carries.
$ws.Range("A1").Value = "aucune instruction"

# Make the new sheet the active tab/sheet (mirrors activeTab 3 -> 4 and
# tabSelected moving off of "2012" and onto "2013").
$ws.Activate() | Out-Null
$ws.Range("E10").Select() | Out-Null
